$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 1496.6666
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 1496.6666
$ws.Range("K32").Value = 0
$ws.Range("L32").ClearContents()
$ws.Range("M32").Value = 1496.6666
$ws.Range("N32").Value = -2148.6666
$ws.Range("H54").Value = 10000
$ws.Range("I54").Value = 0
$ws.Range("K54").Value = 0
$ws.Range("M54").ClearContents()
$ws.Range("H61").Value = 543.125
$ws.Range("J61").Value = 0
$ws.Range("L61").Value = 0
$ws.Range("N61").ClearContents()
$ws.Range("H107").Value = 769.5357
$ws.Range("I107").Value = 568.7917
$ws.Range("J107").Value = 1974
$ws.Range("K107").Value = 568.7917
$ws.Range("L107").Value = 1974
$ws.Range("M107").Value = 1351.2083
$ws.Range("N107").Value = -5814
$ws.Range("H121").Value = 10499.583
$ws.Range("J121").Value = 12439.5
$ws.Range("L121").Value = 37318.5
$ws.Range("N121").Value = -40812.5
$ws.Range("H129").Value = 186222.31
$ws.Range("I129").Value = 298.5
$ws.Range("J129").Value = 201096.22
$ws.Range("K129").Value = 895.5
$ws.Range("L129").Value = 603288.66
$ws.Range("M129").Value = 4104.5
$ws.Range("N129").Value = -613288.66
$ws.Range("H134").Value = 41333
$ws.Range("J134").Value = 41333
$ws.Range("L134").Value = 41333
$ws.Range("N134").Value = -51473
$ws.Range("H138").Value = 28574220
$ws.Range("I138").Value = 66668484
$ws.Range("J138").Value = 3521.35
$ws.Range("K138").Value = 200005452
$ws.Range("L138").Value = 10564.05
$ws.Range("M138").Value = -200000312
$ws.Range("N138").Value = -20844.05
$ws.Range("H140").Value = 50659
$ws.Range("J140").Value = 50659
$ws.Range("L140").Value = 50659
$ws.Range("N140").Value = -61019
$ws.Range("H141").Value = 3287.1667
$ws.Range("I141").Value = 2843.7693
$ws.Range("J141").Value = 4440
$ws.Range("K141").Value = 8531.3079
$ws.Range("L141").Value = 13320
$ws.Range("M141").Value = -3351.3079
$ws.Range("N141").Value = -23680

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 300
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 300
$ws.Range("K4").Value = 0
$ws.Range("L4").ClearContents()
$ws.Range("M4").Value = 300
$ws.Range("N4").Value = -532
$ws.Range("H5").Value = 126.25
$ws.Range("I5").Value = 152.5
$ws.Range("K5").Value = 152.5
$ws.Range("M5").Value = -40.5
$ws.Range("H6").Value = 0
$ws.Range("I6").Value = 0
$ws.Range("K6").Value = 0
$ws.Range("M6").ClearContents()
$ws.Range("H32").Value = 4894.805
$ws.Range("I32").Value = 4148.3716
$ws.Range("K32").Value = 4148.3716
$ws.Range("M32").Value = -3861.3716
$ws.Range("H132").Value = 18288.322
$ws.Range("I132").Value = 1837.04
$ws.Range("J132").Value = 86835.336
$ws.Range("K132").Value = 5511.12
$ws.Range("L132").Value = 260506.008
$ws.Range("M132").Value = -2981.12
$ws.Range("N132").Value = -265566.008

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 126.25
$ws.Range("I4").Value = 152.5
$ws.Range("K4").Value = 152.5
$ws.Range("M4").Value = -37.5
$ws.Range("H11").Value = 5999.5
$ws.Range("I11").Value = 5999.5
$ws.Range("K11").Value = 5999.5
$ws.Range("M11").Value = -5859.5
$ws.Range("H99").Value = 1619.8
$ws.Range("I99").Value = 1649.75
$ws.Range("K99").Value = 1649.75
$ws.Range("M99").Value = -151.75
$ws.Range("H134").Value = 2448.7166
$ws.Range("I134").Value = 2380.4182
$ws.Range("K134").Value = 7141.2546
$ws.Range("M134").Value = -4606.2546

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 27
$ws.Range("J7").Value = 34.666668
$ws.Range("L7").Value = 34.666668
$ws.Range("N7").Value = -260.666668
$ws.Range("H22").Value = 525.7778
$ws.Range("I22").Value = 232.5
$ws.Range("J22").Value = 760.4
$ws.Range("K22").Value = 232.5
$ws.Range("L22").Value = 760.4
$ws.Range("M22").Value = 117.5
$ws.Range("N22").Value = -1460.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H29").Value = 400
$ws.Range("I29").Value = 0
$ws.Range("J29").Value = 400
$ws.Range("K29").Value = 0
$ws.Range("L29").ClearContents()
$ws.Range("M29").Value = 1200
$ws.Range("N29").Value = -1754
$ws.Range("H54").Value = 3753
$ws.Range("I54").Value = 1999
$ws.Range("K54").Value = 5997
$ws.Range("M54").Value = -5438
$ws.Range("H131").Value = 111915.21
$ws.Range("J131").Value = 115743.2
$ws.Range("L131").Value = 347229.6
$ws.Range("N131").Value = -357309.6
$ws.Range("H132").Value = 435
$ws.Range("I132").Value = 400
$ws.Range("K132").Value = 3600
$ws.Range("M132").Value = -1070
$ws.Range("H137").Value = 22228828
$ws.Range("J137").Value = 25648454
$ws.Range("L137").Value = 76945362
$ws.Range("N137").Value = -76955562

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 51.8125
$ws.Range("I2").Value = 49.545456
$ws.Range("J2").Value = 56.8
$ws.Range("K2").Value = 49.545456
$ws.Range("L2").Value = 56.8
$ws.Range("M2").Value = 63.454544
$ws.Range("N2").Value = -282.8
$ws.Range("H11").Value = 12385462
$ws.Range("I11").Value = 13909091
$ws.Range("J11").Value = 4005502
$ws.Range("K11").Value = 13909091
$ws.Range("L11").Value = 4005502
$ws.Range("M11").Value = -13908952
$ws.Range("N11").Value = -4005780
$ws.Range("H97").Value = 1895.5555
$ws.Range("I97").Value = 1765.7142
$ws.Range("J97").Value = 2350
$ws.Range("K97").Value = 1765.7142
$ws.Range("L97").Value = 2350
$ws.Range("M97").Value = -1269.7142
$ws.Range("N97").Value = -3342
$ws.Range("H113").Value = 13101.571
$ws.Range("I113").Value = 17542.2
$ws.Range("J113").Value = 2000
$ws.Range("K113").Value = 17542.2
$ws.Range("L113").Value = 2000
$ws.Range("M113").Value = -15372.2
$ws.Range("N113").Value = -6340
$ws.Range("H126").Value = 5142.273
$ws.Range("I126").Value = 4087.2273
$ws.Range("J126").Value = 7252.364
$ws.Range("K126").Value = 12261.6819
$ws.Range("L126").Value = 21757.092
$ws.Range("M126").Value = -9791.6819
$ws.Range("N126").Value = -26697.092
$ws.Range("H141").Value = 44859.6
$ws.Range("J141").Value = 44859.6
$ws.Range("L141").Value = 44859.6
$ws.Range("N141").Value = -55219.6

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3010.3447
$ws.Range("I40").Value = 2485.3
$ws.Range("K40").Value = 2485.3
$ws.Range("M40").Value = -2349.3
$ws.Range("H46").Value = 1034.2433
$ws.Range("I46").Value = 1010.5143
$ws.Range("J46").Value = 1449.5
$ws.Range("K46").Value = 1010.5143
$ws.Range("L46").Value = 1449.5
$ws.Range("M46").Value = -822.5143
$ws.Range("N46").Value = -1825.5
$ws.Range("H122").Value = 656277.3
$ws.Range("I122").Value = 1227740.2
$ws.Range("J122").Value = 3176.7856
$ws.Range("K122").Value = 3683220.6
$ws.Range("L122").Value = 9530.356800000001
$ws.Range("M122").Value = -3680770.6
$ws.Range("N122").Value = -14430.3568

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H24").Value = 10000000
$ws.Range("I24").Value = 0
$ws.Range("K24").Value = 0
$ws.Range("M24").ClearContents()
$ws.Range("H54").Value = 17500
$ws.Range("J54").Value = 17500
$ws.Range("L54").Value = 17500
$ws.Range("N54").Value = -18540
$ws.Range("H126").Value = 1230.125
$ws.Range("I126").Value = 1295.2142
$ws.Range("J126").Value = 774.5
$ws.Range("K126").Value = 3885.6426
$ws.Range("L126").Value = 2323.5
$ws.Range("M126").Value = -1415.6426
$ws.Range("N126").Value = -7263.5
$ws.Range("H141").Value = 63333.332
$ws.Range("J141").Value = 63333.332
$ws.Range("L141").Value = 63333.332
$ws.Range("N141").Value = -73693.33199999999
